# Update cryptocurrency price/volume figures to the latest scrape.
# Cells are plain text (e.g. "328.15", "-1.00%") in the source workbook, so we
# temporarily mark each cell as Text before assigning, then restore its original
# style -- this keeps the values as literal strings (matching the source data)
# instead of letting Excel auto-convert them into numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '328.15' }
    @{ Cell = 'E2'; Value = '-1.00%' }
    @{ Cell = 'D3'; Value = '44.31' }
    @{ Cell = 'E3'; Value = '-0.82%' }
    @{ Cell = 'D4'; Value = '5.347' }
    @{ Cell = 'E4'; Value = '-3.75%' }
    @{ Cell = 'D5'; Value = '0.08368' }
    @{ Cell = 'E5'; Value = '1.09%' }
    @{ Cell = 'D6'; Value = '1.940' }
    @{ Cell = 'E6'; Value = '-4.87%' }
    @{ Cell = 'D7'; Value = '0.9722' }
    @{ Cell = 'E7'; Value = '-0.50%' }
    @{ Cell = 'D9'; Value = '0.1123' }
    @{ Cell = 'E9'; Value = '-0.18%' }
    @{ Cell = 'D10'; Value = '0.1906' }
    @{ Cell = 'E10'; Value = '0.29%' }
    @{ Cell = 'D11'; Value = '0.09684' }
    @{ Cell = 'E11'; Value = '-3.04%' }
    @{ Cell = 'D12'; Value = '0.04598' }
    @{ Cell = 'E12'; Value = '-1.50%' }
    @{ Cell = 'E13'; Value = '0.37%' }
    @{ Cell = 'D14'; Value = '0.001302' }
    @{ Cell = 'E14'; Value = '2.28%' }
    @{ Cell = 'D15'; Value = '0.005839' }
    @{ Cell = 'E15'; Value = '-1.31%' }
    @{ Cell = 'D16'; Value = '3.360' }
    @{ Cell = 'E16'; Value = '-0.05%' }
    @{ Cell = 'D17'; Value = '4.417' }
    @{ Cell = 'E17'; Value = '-0.41%' }
    @{ Cell = 'E18'; Value = '0.17%' }
    @{ Cell = 'D19'; Value = '8.301' }
    @{ Cell = 'E19'; Value = '-19.33%' }
    @{ Cell = 'E20'; Value = '0.26%' }
    @{ Cell = 'D21'; Value = '0.2654' }
    @{ Cell = 'E21'; Value = '6.50%' }
    @{ Cell = 'D22'; Value = '0.04181' }
    @{ Cell = 'D23'; Value = '0.001242' }
    @{ Cell = 'E23'; Value = '-4.58%' }
    @{ Cell = 'D24'; Value = '0.004413' }
    @{ Cell = 'E24'; Value = '0.06%' }
    @{ Cell = 'E25'; Value = '1.43%' }
    @{ Cell = 'D26'; Value = '0.0002979' }
    @{ Cell = 'E26'; Value = '-20.36%' }
    @{ Cell = 'D38'; Value = '0.02720' }
    @{ Cell = 'E38'; Value = '-2.56%' }
    @{ Cell = 'D39'; Value = '0.05639' }
    @{ Cell = 'E39'; Value = '-1.90%' }
    @{ Cell = 'D40'; Value = '0.007822' }
    @{ Cell = 'E40'; Value = '2.65%' }
    @{ Cell = 'D41'; Value = '0.1413' }
    @{ Cell = 'E41'; Value = '-0.73%' }
    @{ Cell = 'D42'; Value = '0.007348' }
    @{ Cell = 'E42'; Value = '-2.79%' }
    @{ Cell = 'E43'; Value = '3.27%' }
    @{ Cell = 'D44'; Value = '0.008704' }
    @{ Cell = 'E44'; Value = '4.67%' }
    @{ Cell = 'D45'; Value = '0.3513' }
    @{ Cell = 'D46'; Value = '0.00006909' }
    @{ Cell = 'E46'; Value = '-1.91%' }
    @{ Cell = 'E47'; Value = '-0.10%' }
    @{ Cell = 'D48'; Value = '0.003485' }
    @{ Cell = 'E48'; Value = '-2.75%' }
    @{ Cell = 'D49'; Value = '0.003531' }
    @{ Cell = 'E49'; Value = '39.89%' }
    @{ Cell = 'D50'; Value = '0.00002101' }
    @{ Cell = 'E50'; Value = '-0.10%' }
    @{ Cell = 'D51'; Value = '0.0002001' }
    @{ Cell = 'E51'; Value = '-0.10%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
    $rng.Style = $origStyle
}

